# Update column F ("dSF") values per repulled data.
# Mapping of row number -> new value for column F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    4  = -2
    5  = 4
    6  = -1
    7  = 1
    8  = 6
    9  = -1
    10 = -1
    12 = -2
    13 = 0
    15 = -2
    16 = 1
    17 = -1
    18 = 3
    19 = 4
    20 = 2
    21 = -2
    22 = 4
    23 = 2
    25 = 6
    26 = -1
    27 = -3
    28 = -3
    30 = 5
    31 = 8
    32 = 4
    33 = 1
    34 = 1
    35 = -1
    36 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
